# Apply the edit described by the diff:
# A new data row is inserted at row 37 (shifting all following rows down by one,
# the last existing row 151 becomes row 152), and the new row 37 is populated
# with a new Espinaca price record dated 2021-09-20 (serial 44459).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37; Excel shifts rows 37..151 down to 38..152
# and carries formatting down from the row above (so D37 keeps the date style).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's values.
$ws.Cells.Item(37, 1).Value = 8
$ws.Cells.Item(37, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = 44459
$ws.Cells.Item(37, 5).Value = 4
$ws.Cells.Item(37, 6).Value = 100112012
$ws.Cells.Item(37, 7).Value = "Espinaca"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 2000
$ws.Cells.Item(37, 11).Value = 400
$ws.Cells.Item(37, 12).Value = 500
$ws.Cells.Item(37, 13).Value = 450
$ws.Cells.Item(37, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(37, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 16).Value = 900
$ws.Cells.Item(37, 17).Value = 0.5
$ws.Cells.Item(37, 18).Value = "Hortaliza"
